# Fix Uganda life-expectancy "Data" sheet: populate the blank "entity-name"
# column (B) with the district name that corresponds to each row's district
# id (A), using the official district-id -> district-name lookup table.

$districtNames = @{
    "d101" = "Kalangala"
    "d102" = "Kampala"
    "d103" = "Kiboga"
    "d104" = "Luwero"
    "d105" = "Masaka"
    "d106" = "Mpigi"
    "d107" = "Mubende"
    "d108" = "Mukono"
    "d109" = "Nakasongola"
    "d110" = "Rakai"
    "d111" = "Ssembabule"
    "d112" = "Kayunga"
    "d113" = "Wakiso"
    "d114" = "Lyantonde"
    "d115" = "Mityana"
    "d116" = "Nakaseke"
    "d117" = "Buikwe"
    "d118" = "Bukomansimbi"
    "d119" = "Butambala"
    "d120" = "Buvuma"
    "d121" = "Gomba"
    "d122" = "Kalungu"
    "d123" = "Kyankwanzi"
    "d124" = "Lwengo"
    "d201" = "Bugiri"
    "d202" = "Busia"
    "d203" = "Iganga"
    "d204" = "Jinja"
    "d205" = "Kamuli"
    "d206" = "Kapchorwa"
    "d207" = "Katakwi"
    "d208" = "Kumi"
    "d209" = "Mbale"
    "d210" = "Pallisa"
    "d211" = "Soroti"
    "d212" = "Tororo"
    "d213" = "Kaberamaido"
    "d214" = "Mayuge"
    "d215" = "Sironko"
    "d216" = "Amuria"
    "d217" = "Budaka"
    "d218" = "Bududa"
    "d219" = "Bukedea"
    "d220" = "Bukwo"
    "d221" = "Butaleja"
    "d222" = "Kaliro"
    "d223" = "Manafwa"
    "d224" = "Namutumba"
    "d225" = "Bulambuli"
    "d226" = "Buyende"
    "d227" = "Kibuku"
    "d228" = "Kween"
    "d229" = "Luuka"
    "d230" = "Namayingo"
    "d231" = "Ngora"
    "d232" = "Serere"
    "d301" = "Adjumani"
    "d302" = "Apac"
    "d303" = "Arua"
    "d304" = "Gulu"
    "d305" = "Kitgum"
    "d306" = "Kotido"
    "d307" = "Lira"
    "d308" = "Moroto"
    "d309" = "Moyo"
    "d310" = "Nebbi"
    "d311" = "Nakapiripirit"
    "d312" = "Pader"
    "d313" = "Yumbe"
    "d314" = "Abim"
    "d315" = "Amolatar"
    "d316" = "Amuru"
    "d317" = "Dokolo"
    "d318" = "Kaabong"
    "d319" = "Koboko"
    "d320" = "Maracha"
    "d321" = "Oyam"
    "d322" = "Agago"
    "d323" = "Alebtong"
    "d324" = "Amudat"
    "d325" = "Kole"
    "d326" = "Lamwo"
    "d327" = "Napak"
    "d328" = "Nwoya"
    "d329" = "Otuke"
    "d330" = "Zombo"
    "d401" = "Bundibugyo"
    "d402" = "Bushenyi"
    "d403" = "Hoima"
    "d404" = "Kabale"
    "d405" = "Kabarole"
    "d406" = "Kasese"
    "d407" = "Kibaale"
    "d408" = "Kisoro"
    "d409" = "Masindi"
    "d410" = "Mbarara"
    "d411" = "Ntungamo"
    "d412" = "Rukungiri"
    "d413" = "Kamwenge"
    "d414" = "Kanungu"
    "d415" = "Kyenjojo"
    "d416" = "Buliisa"
    "d417" = "Ibanda"
    "d418" = "Isingiro"
    "d419" = "Kiruhura"
    "d420" = "Buhweju"
    "d421" = "Kiryandongo"
    "d422" = "Kyegegwa"
    "d423" = "Mitooma"
    "d424" = "Ntoroko"
    "d425" = "Rubirizi"
    "d426" = "Sheema"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $id = $ws.Cells.Item($r, 1).Text
    if ($districtNames.ContainsKey($id)) {
        $ws.Cells.Item($r, 2).Value = $districtNames[$id]
    }
}
